$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1494
$ws.Range("F5").Value = 7521
$ws.Range("F7").Value = 4800
$ws.Range("F8").Value = 7062
$ws.Range("F10").Value = 284
$ws.Range("F11").Value = 1496
$ws.Range("F13").Value = 182
$ws.Range("F15").Value = 1168
$ws.Range("F17").Value = 171
$ws.Range("F22").Value = 1172
$ws.Range("F23").Value = 954
$ws.Range("F26").Value = 1231
$ws.Range("F27").Value = 47
$ws.Range("F31").Value = 186
$ws.Range("F34").Value = 104
$ws.Range("F37").Value = 426
$ws.Range("F39").Value = 64
$ws.Range("F40").Value = 383
$ws.Range("F41").Value = 1202

# Sheet "演出" (Shows)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F27").Value = 638
$ws.Range("F29").Value = 31
$ws.Range("F33").Value = 990
$ws.Range("F39").Value = 13

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 855
$ws.Range("F6").Value = 678
$ws.Range("F8").Value = 1621
$ws.Range("F9").Value = 2524

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1494
$ws.Range("F4").Value = 855
$ws.Range("F7").Value = 678
$ws.Range("F8").Value = 678
$ws.Range("F9").Value = 7521
$ws.Range("F11").Value = 4800
$ws.Range("F13").Value = 7062
$ws.Range("F14").Value = 284
$ws.Range("F15").Value = 1496
$ws.Range("F17").Value = 182
$ws.Range("F18").Value = 1621
$ws.Range("F19").Value = 2524
$ws.Range("F22").Value = 1168
$ws.Range("F23").Value = 171
$ws.Range("F26").Value = 1172
$ws.Range("F27").Value = 638
$ws.Range("F28").Value = 954
$ws.Range("F30").Value = 1231
$ws.Range("F32").Value = 186
$ws.Range("F33").Value = 31
$ws.Range("F36").Value = 104
$ws.Range("F37").Value = 990
$ws.Range("F41").Value = 64
$ws.Range("F43").Value = 383
